# Remove the blank paragraph and the "Please read Using Pressure Canners..."
# paragraph that follow the "Quality: ..." paragraph, collapsing the
# document back down to a single trailing blank paragraph after "Quality".

$d = $word.ActiveDocument

$target = "Please read Using Pressure Canners before beginning."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$target*") {
        # The paragraph immediately before this one is the blank "\n"
        # paragraph that also needs to go.
        $blank = $d.Paragraphs.Item($i - 1)
        $range = $d.Range($blank.Range.Start, $p.Range.End)
        $range.Delete()
        break
    }
}
